$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("In Advance" data shifts from N->O,
# "Outstanding" data shifts from P->Q), mirroring the author's "Variable Instalments"
# column insert for the RBI loan schedule.
$ws.Columns("N:N").Insert()

# Match the width Excel copies in from the column immediately to the left (M) when
# inserting a column via the UI.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment Schedule" the active sheet (was "Transactions") and update its
# selected cell.
[void]$ws.Activate()
[void]$ws.Range("J20").Select()
